$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "Sistemas"
$ws.Range("D3").Value = "Sistemas"

$ws.Range("E6").Select()
